$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header / summary cells ---
# VALOR MORA total (row 11)
$ws.Cells.Item(11, 5).Value2 = 394908

# Cant. Trabajadores (row13 C) and Cant. Periodos (row13 F)
$ws.Cells.Item(13, 3).Value2 = 2
$ws.Cells.Item(13, 6).Value2 = 9

# --- Remove the 7 surplus worker/period detail rows ---
# Deleting rows 16:22 shifts the remaining detail rows (old 23-34) up so
# that the final detail row keeps its special "closing" border style, and
# the signature block (old rows 39-40) moves up to rows 32-33, matching
# the new layout.
$ws.Rows("16:22").Delete()

# --- Overwrite the detail table (now rows 16-27) with the new contents ---
$data = @(
    @("PE", "900063325041982", "DARWIN ALBERTO FERNANDEZ TORRES", "1809", 32909, 877803),
    @("PE", "900063325041982", "DARWIN ALBERTO FERNANDEZ TORRES", "1810", 32909, 877803),
    @("PE", "900063325041982", "DARWIN ALBERTO FERNANDEZ TORRES", "1811", 32909, 877803),
    @("PE", "900063325041982", "DARWIN ALBERTO FERNANDEZ TORRES", "1812", 32909, 877803),
    @("PE", "900063325041982", "DARWIN ALBERTO FERNANDEZ TORRES", "1901", 32909, 877803),
    @("PE", "941236224021987", "JOHAN JOSE VILLALOBOS HERNANDEZ", "1902", 32909, 877803),
    @("PE", "900063325041982", "DARWIN ALBERTO FERNANDEZ TORRES", "1902", 32909, 877803),
    @("PE", "941236224021987", "JOHAN JOSE VILLALOBOS HERNANDEZ", "1903", 32909, 877803),
    @("PE", "900063325041982", "DARWIN ALBERTO FERNANDEZ TORRES", "1903", 32909, 877803),
    @("PE", "941236224021987", "JOHAN JOSE VILLALOBOS HERNANDEZ", "1904", 32909, 877803),
    @("PE", "900063325041982", "DARWIN ALBERTO FERNANDEZ TORRES", "1904", 32909, 877803),
    @("PE", "900063325041982", "DARWIN ALBERTO FERNANDEZ TORRES", "1905", 32909, 877803)
)

$r = 16
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value2 = $row[0]
    $ws.Cells.Item($r, 3).Value2 = $row[1]
    $ws.Cells.Item($r, 4).Value2 = $row[2]
    $ws.Cells.Item($r, 5).Value2 = $row[3]
    $ws.Cells.Item($r, 6).Value2 = $row[4]
    $ws.Cells.Item($r, 7).Value2 = $row[5]
    $r = $r + 1
}
